# Generate Report for Handoff
# Adds a second data row (for file "ffffc3757e19-ef34-4fa2-9600-110385f3efa9.md")
# to each of the three tables (Overview / zh-cn / de-de) alongside the
# existing row for "f2cd8a88-8020-4b70-ab5b-5a958ea4bffe.md" (renamed from
# the old "d663a4a2-1e0e-46ed-8636-3522b2a85a31.md" handoff).

$wb = $excel.ActiveWorkbook

$oldBase = "d663a4a2-1e0e-46ed-8636-3522b2a85a31"
$newBase = "f2cd8a88-8020-4b70-ab5b-5a958ea4bffe"
$newHash = "6a42664e605b5a09a9b0fb1f6178ead5077864f9"
$dupBase = "ffffc3757e19-ef34-4fa2-9600-110385f3efa9"

$handoffDate = "2016-09-09 12:47:47"
$xliffDate   = "2016-09-09 12:47:36"

# ---------------------------------------------------------------------
# 1. Rename the existing handed-off file throughout the workbook
#    (d663a4a2-... -> f2cd8a88-...), and bump its timestamps.
# ---------------------------------------------------------------------

# -- Overview sheet --
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newBase.md"
$wsOverview.Hyperlinks.Item(1).Address = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6209404a21d795e4215075db29fe7dc4c7210b98/e2e/$newBase.md"
$wsOverview.Hyperlinks.Item(1).TextToDisplay = "e2e\$newBase.md"
$wsOverview.Range("G2").Value = $handoffDate

# -- zh-cn sheet --
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = "$newBase.md"
$wsZh.Hyperlinks.Item(1).Address = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6209404a21d795e4215075db29fe7dc4c7210b98/e2e/$newBase.md"
$wsZh.Hyperlinks.Item(1).TextToDisplay = "$newBase.md"
$wsZh.Range("G2").Value = "$newBase.$newHash.zh-cn.xlf"
$wsZh.Range("H2").Value = $xliffDate

# -- de-de sheet --
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = "$newBase.md"
$wsDe.Hyperlinks.Item(1).Address = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6209404a21d795e4215075db29fe7dc4c7210b98/e2e/$newBase.md"
$wsDe.Hyperlinks.Item(1).TextToDisplay = "$newBase.md"
$wsDe.Range("G2").Value = "$newBase.$newHash.de-de.xlf"
$wsDe.Range("H2").Value = $handoffDate

# ---------------------------------------------------------------------
# 2. Append a new row (duplicate handoff) to each table for the
#    "ffffc3757e19-..." file.
# ---------------------------------------------------------------------

# -- Overview: File Name | Path And Name | Extension | Publish URL | zh-cn | de-de | Latest HO Xliff Generate Date --
$loOverview = $wsOverview.ListObjects.Item(1)
$rowOverview = $loOverview.ListRows.Add()
$wsOverview.Range("A3").Value = "$dupBase.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = $handoffDate
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6209404a21d795e4215075db29fe7dc4c7210b98/e2e/$dupBase.md", "", "", "e2e\$dupBase.md") | Out-Null

# -- zh-cn / de-de: Source File Name | File Extension | Status | Source Path | Priority |
#    Content Duplicate | Latest Handoff File | Latest Handoff Datetime | Latest Target File |
#    Latest Handback File | Latest Handback DateTime | Reference Tokens | To be localized |
#    Dependency From | Has metadata | Error Detail --
$loZh = $wsZh.ListObjects.Item(1)
$rowZh = $loZh.ListRows.Add()
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "True"
$wsZh.Range("G3").Value = "$newBase.$newHash.zh-cn.xlf"
$wsZh.Range("H3").Value = $xliffDate
$wsZh.Range("K3").Value = "0001-01-01 00:00:00"
$wsZh.Range("M3").Value = "True"
$wsZh.Range("O3").Value = "False"
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6209404a21d795e4215075db29fe7dc4c7210b98/e2e/$dupBase.md", "", "", "$dupBase.md") | Out-Null

$loDe = $wsDe.ListObjects.Item(1)
$rowDe = $loDe.ListRows.Add()
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "True"
$wsDe.Range("G3").Value = "$newBase.$newHash.de-de.xlf"
$wsDe.Range("H3").Value = $handoffDate
$wsDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDe.Range("M3").Value = "True"
$wsDe.Range("O3").Value = "False"
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6209404a21d795e4215075db29fe7dc4c7210b98/e2e/$dupBase.md", "", "", "$dupBase.md") | Out-Null
